$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.192.51"
$ws.Range("E2").Value = "  -0.31%  "

# Row 3
$ws.Range("D3").Value = "1.815.46"
$ws.Range("E3").Value = "  +1.69%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.87"
$ws.Range("E5").Value = "  +0.18%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.557"
$ws.Range("E6").Value = "  +1.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.26%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.10"
$ws.Range("E8").Value = "  -4.62%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.292"
$ws.Range("E9").Value = "  +3.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0742"
$ws.Range("E10").Value = "  +12.91%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0929"
$ws.Range("E11").Value = "  -0.23%  "

# Row 12
$ws.Range("D12").Value = "2.074.71"
$ws.Range("E12").Value = "  +1.61%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.12"
$ws.Range("E13").Value = "  +0.06%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.815.51"
$ws.Range("E14").Value = "  +1.66%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.643"
$ws.Range("E15").Value = "  +1.68%  "

# Row 16
$ws.Range("D16").Value = "34.203.07"
$ws.Range("E16").Value = "  -0.22%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.33"
$ws.Range("E17").Value = "  +2.43%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.68"
$ws.Range("E18").Value = "  +1.14%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "250.03"
$ws.Range("E19").Value = "  -1.95%  "

# Row 20
$ws.Range("E20").Value = "  +9.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.10"
$ws.Range("E21").Value = "  +6.43%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.28"
$ws.Range("E23").Value = "  +2.45%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.16"
$ws.Range("E24").Value = "  +1.50%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.47"
$ws.Range("E25").Value = "  +1.72%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.70"
$ws.Range("E26").Value = "  +1.76%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.25"
$ws.Range("E27").Value = "  +3.36%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.115"
$ws.Range("E28").Value = "  +0.93%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.24%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0533"
$ws.Range("E30").Value = "  +3.65%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.79"
$ws.Range("E31").Value = "  +0.32%  "

# Row 32
$ws.Range("E32").Value = "  +1.96%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.60"
$ws.Range("E33").Value = "  +0.53%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.89"
$ws.Range("E34").Value = "  +0.06%  "

# Row 35
$ws.Range("D35").Value = "1.432.97"
$ws.Range("E35").Value = "  -1.01%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.643"
$ws.Range("E36").Value = "  +3.17%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.06"
$ws.Range("E37").Value = "  +1.08%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0191"
$ws.Range("E38").Value = "  +1.57%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.963"
$ws.Range("E39").Value = "  +8.16%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "81.83"
$ws.Range("E40").Value = "  -1.55%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.75"
$ws.Range("E41").Value = "  -3.46%  "

# Row 42
$ws.Range("E42").Value = "  -0.42%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.18"
$ws.Range("E43").Value = "  +4.63%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.07"
$ws.Range("E44").Value = "  +2.28%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0498"
$ws.Range("E45").Value = "  -1.88%  "

# Row 46
$ws.Range("D46").Value = "1.971.01"
$ws.Range("E46").Value = "  +1.44%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.24"
$ws.Range("E47").Value = "  +7.80%  "

# Row 48
$ws.Range("E48").Value = "  -1.42%  "

# Row 49
$ws.Range("E49").Value = "  -0.40%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.85"
$ws.Range("E50").Value = "  -4.11%  "

# Row 51
$ws.Range("D51").Value = "0.0₆0124"
$ws.Range("E51").Value = "  +5.90%  "
